$wb = $excel.ActiveWorkbook

# The second localization file (95c62293-...) has finished translation and is
# now ready to be handed back, for both the zh-cn and de-de locales. Update
# the per-locale status/handoff-datetime rows and roll the summary up into
# the Overview sheet.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-15-20 08:15:51"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-20 08:15:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-20 08:15:51"
